# Auto-generated script to update 'Recorded By' (column G) values
# per the target diff: reorder comma-separated entries in specific rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Row 2: 'System, system, backup@backdoor.com' -> 'System, backup@backdoor.com, system'
$ws.Cells.Item(2, 7).Value = 'System, backup@backdoor.com, system'

# Row 3: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(3, 7).Value = 'dnasr281@gmail.com, System'

# Row 6: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(6, 7).Value = 'dnasr281@gmail.com, System'

# Row 10: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(10, 7).Value = 'dnasr281@gmail.com, System'

# Row 11: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(11, 7).Value = 'dnasr281@gmail.com, System'

# Row 12: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(12, 7).Value = 'dnasr281@gmail.com, System'

# Row 13: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(13, 7).Value = 'dnasr281@gmail.com, System'

# Row 14: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(14, 7).Value = 'dnasr281@gmail.com, System'

# Row 15: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(15, 7).Value = 'dnasr281@gmail.com, System'

# Row 17: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(17, 7).Value = 'dnasr281@gmail.com, System'

# Row 18: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(18, 7).Value = 'dnasr281@gmail.com, System'

# Row 19: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(19, 7).Value = 'dnasr281@gmail.com, System'

# Row 20: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(20, 7).Value = 'dnasr281@gmail.com, System'

# Row 21: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(21, 7).Value = 'dnasr281@gmail.com, System'

# Row 22: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(22, 7).Value = 'dnasr281@gmail.com, System'

# Row 24: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(24, 7).Value = 'dnasr281@gmail.com, System'

# Row 26: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(26, 7).Value = 'dnasr281@gmail.com, System'

# Row 28: 'System, system, backup@backdoor.com' -> 'System, backup@backdoor.com, system'
$ws.Cells.Item(28, 7).Value = 'System, backup@backdoor.com, system'

# Row 29: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(29, 7).Value = 'dnasr281@gmail.com, System'

# Row 32: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(32, 7).Value = 'dnasr281@gmail.com, System'

# Row 36: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(36, 7).Value = 'dnasr281@gmail.com, System'

# Row 37: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(37, 7).Value = 'dnasr281@gmail.com, System'

# Row 38: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(38, 7).Value = 'dnasr281@gmail.com, System'

# Row 39: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(39, 7).Value = 'dnasr281@gmail.com, System'

# Row 40: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(40, 7).Value = 'dnasr281@gmail.com, System'

# Row 41: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(41, 7).Value = 'dnasr281@gmail.com, System'

# Row 43: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(43, 7).Value = 'dnasr281@gmail.com, System'

# Row 44: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(44, 7).Value = 'dnasr281@gmail.com, System'

# Row 45: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(45, 7).Value = 'dnasr281@gmail.com, System'

# Row 46: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(46, 7).Value = 'dnasr281@gmail.com, System'

# Row 47: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(47, 7).Value = 'dnasr281@gmail.com, System'

# Row 48: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(48, 7).Value = 'dnasr281@gmail.com, System'

# Row 50: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(50, 7).Value = 'dnasr281@gmail.com, System'

# Row 52: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(52, 7).Value = 'dnasr281@gmail.com, System'

# Row 54: 'System, system, backup@backdoor.com' -> 'System, backup@backdoor.com, system'
$ws.Cells.Item(54, 7).Value = 'System, backup@backdoor.com, system'

# Row 55: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(55, 7).Value = 'dnasr281@gmail.com, System'

# Row 58: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(58, 7).Value = 'dnasr281@gmail.com, System'

# Row 62: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(62, 7).Value = 'dnasr281@gmail.com, System'

# Row 63: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(63, 7).Value = 'dnasr281@gmail.com, System'

# Row 64: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(64, 7).Value = 'dnasr281@gmail.com, System'

# Row 65: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(65, 7).Value = 'dnasr281@gmail.com, System'

# Row 66: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(66, 7).Value = 'dnasr281@gmail.com, System'

# Row 67: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(67, 7).Value = 'dnasr281@gmail.com, System'

# Row 69: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(69, 7).Value = 'dnasr281@gmail.com, System'

# Row 70: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(70, 7).Value = 'dnasr281@gmail.com, System'

# Row 71: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(71, 7).Value = 'dnasr281@gmail.com, System'

# Row 72: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(72, 7).Value = 'dnasr281@gmail.com, System'

# Row 73: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(73, 7).Value = 'dnasr281@gmail.com, System'

# Row 74: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(74, 7).Value = 'dnasr281@gmail.com, System'

# Row 76: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(76, 7).Value = 'dnasr281@gmail.com, System'

# Row 78: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(78, 7).Value = 'dnasr281@gmail.com, System'

# Row 83: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(83, 7).Value = 'dnasr281@gmail.com, System'

# Row 84: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(84, 7).Value = 'dnasr281@gmail.com, System'

# Row 85: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(85, 7).Value = 'dnasr281@gmail.com, System'

# Row 86: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(86, 7).Value = 'dnasr281@gmail.com, System'

# Row 87: 'admin@admin.com, dnasr281@gmail.com' -> 'dnasr281@gmail.com, admin@admin.com'
$ws.Cells.Item(87, 7).Value = 'dnasr281@gmail.com, admin@admin.com'

# Row 90: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(90, 7).Value = 'dnasr281@gmail.com, System'

# Row 92: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(92, 7).Value = 'dnasr281@gmail.com, System'

# Row 93: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(93, 7).Value = 'dnasr281@gmail.com, System'

# Row 94: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(94, 7).Value = 'dnasr281@gmail.com, System'

# Row 96: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(96, 7).Value = 'dnasr281@gmail.com, System'

# Row 99: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(99, 7).Value = 'dnasr281@gmail.com, System'

# Row 101: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(101, 7).Value = 'dnasr281@gmail.com, System'

# Row 109: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(109, 7).Value = 'dnasr281@gmail.com, System'

# Row 110: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(110, 7).Value = 'dnasr281@gmail.com, System'

# Row 111: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(111, 7).Value = 'dnasr281@gmail.com, System'

# Row 112: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(112, 7).Value = 'dnasr281@gmail.com, System'

# Row 113: 'admin@admin.com, dnasr281@gmail.com' -> 'dnasr281@gmail.com, admin@admin.com'
$ws.Cells.Item(113, 7).Value = 'dnasr281@gmail.com, admin@admin.com'

# Row 116: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(116, 7).Value = 'dnasr281@gmail.com, System'

# Row 118: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(118, 7).Value = 'dnasr281@gmail.com, System'

# Row 119: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(119, 7).Value = 'dnasr281@gmail.com, System'

# Row 120: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(120, 7).Value = 'dnasr281@gmail.com, System'

# Row 122: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(122, 7).Value = 'dnasr281@gmail.com, System'

# Row 125: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(125, 7).Value = 'dnasr281@gmail.com, System'

# Row 127: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(127, 7).Value = 'dnasr281@gmail.com, System'

# Row 135: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(135, 7).Value = 'dnasr281@gmail.com, System'

# Row 136: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(136, 7).Value = 'dnasr281@gmail.com, System'

# Row 137: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(137, 7).Value = 'dnasr281@gmail.com, System'

# Row 138: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(138, 7).Value = 'dnasr281@gmail.com, System'

# Row 139: 'admin@admin.com, dnasr281@gmail.com' -> 'dnasr281@gmail.com, admin@admin.com'
$ws.Cells.Item(139, 7).Value = 'dnasr281@gmail.com, admin@admin.com'

# Row 142: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(142, 7).Value = 'dnasr281@gmail.com, System'

# Row 144: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(144, 7).Value = 'dnasr281@gmail.com, System'

# Row 145: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(145, 7).Value = 'dnasr281@gmail.com, System'

# Row 146: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(146, 7).Value = 'dnasr281@gmail.com, System'

# Row 148: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(148, 7).Value = 'dnasr281@gmail.com, System'

# Row 151: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(151, 7).Value = 'dnasr281@gmail.com, System'

# Row 153: 'System, dnasr281@gmail.com' -> 'dnasr281@gmail.com, System'
$ws.Cells.Item(153, 7).Value = 'dnasr281@gmail.com, System'
